$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.316.41'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.09%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.648.85'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.20%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.56'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.35%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("E8").Value = '  -0.47%  '

$ws.Range("E9").Value = '  +3.03%  '

$ws.Range("E10").Value = '  -1.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.26'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.68%  '

$ws.Range("E12").Value = '  +1.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.10'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.83%  '

$ws.Range("E14").Value = '  +1.10%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.131.36'
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '68.153.74'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.26%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.632.61'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.26%  '

$ws.Range("E18").Value = '  -0.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '364.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.25%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.37'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.28%  '

$ws.Range("E21").Value = '  +4.11%  '

$ws.Range("E22").Value = '  -0.81%  '

$ws.Range("E23").Value = '  -1.93%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '75.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.62%  '

$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.75'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.36%  '

$ws.Range("E28").Value = '  +0.81%  '

$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '559.10'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.35%  '

$ws.Range("E31").Value = '  +1.04%  '

$ws.Range("E32").Value = '  -1.35%  '

$ws.Range("E33").Value = '  +0.31%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.129'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.34%  '

$ws.Range("E35").Value = '  -0.23%  '

$ws.Range("E36").Value = '  +1.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '160.53'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.84'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.21%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.371'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.04%  '

$ws.Range("E40").Value = '  -2.22%  '

$ws.Range("E41").Value = '  -0.65%  '

$ws.Range("E42").Value = '  +4.81%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.62'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.16%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '159.28'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.43%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.75'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.38%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.15'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.37%  '

$ws.Range("E48").Value = '  -0.71%  '

$ws.Range("E49").Value = '  +0.40%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.613'
$ws.Range("D50").Style = "Normal"

$ws.Range("E51").Value = '  +0.99%  '
